$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 13
$ws.Range("H13").Value = 3000
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
# Row 21
$ws.Range("H21").Value = 7824.875
$ws.Range("I21").Value = 3474.75
$ws.Range("J21").Value = 12175
$ws.Range("K21").Value = 3474.75
$ws.Range("L21").Value = 12175
$ws.Range("M21").Value = -3006.75
$ws.Range("N21").Value = -13111
# Row 23
$ws.Range("H23").Value = 7824.875
$ws.Range("I23").Value = 3474.75
$ws.Range("J23").Value = 12175
$ws.Range("K23").Value = 3474.75
$ws.Range("L23").Value = 12175
$ws.Range("M23").Value = -3240.75
$ws.Range("N23").Value = -12643
# Row 64
$ws.Range("H64").Value = 3168.7878
$ws.Range("I64").Value = 3081.6
$ws.Range("J64").Value = 3184.3572
$ws.Range("K64").Value = 3081.6
$ws.Range("L64").Value = 3184.3572
$ws.Range("M64").Value = -2833.6
$ws.Range("N64").Value = -3680.3572
# Row 67
$ws.Range("H67").Value = 3168.7878
$ws.Range("I67").Value = 3081.6
$ws.Range("J67").Value = 3184.3572
$ws.Range("K67").Value = 3081.6
$ws.Range("L67").Value = 3184.3572
$ws.Range("M67").Value = -2223.6
$ws.Range("N67").Value = -4900.3572
# Row 88
$ws.Range("H88").Value = 7498
$ws.Range("I88").Value = 11612.5
$ws.Range("J88").Value = 3383.5
$ws.Range("K88").Value = 11612.5
$ws.Range("L88").Value = 3383.5
$ws.Range("M88").Value = -11206.5
$ws.Range("N88").Value = -4195.5
# Row 91
$ws.Range("H91").Value = 7498
$ws.Range("I91").Value = 11612.5
$ws.Range("J91").Value = 3383.5
$ws.Range("K91").Value = 11612.5
$ws.Range("L91").Value = 3383.5
$ws.Range("M91").Value = -10208.5
$ws.Range("N91").Value = -6191.5
# Row 137
$ws.Range("H137").Value = 1658.8214
$ws.Range("I137").Value = 1244
$ws.Range("J137").Value = 2534.5557
$ws.Range("K137").Value = 3732
$ws.Range("L137").Value = 7603.6671
$ws.Range("M137").Value = -1182
$ws.Range("N137").Value = -12703.6671
# Row 138
$ws.Range("H138").Value = 417843.16
$ws.Range("I138").Value = 435575.47
$ws.Range("J138").Value = 10000
$ws.Range("K138").Value = 1306726.41
$ws.Range("L138").Value = 30000
$ws.Range("M138").Value = -1301586.41
$ws.Range("N138").Value = -40280

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 9690.796
$ws.Range("I32").Value = 2847.6924
$ws.Range("J32").Value = 29030
$ws.Range("K32").Value = 2847.6924
$ws.Range("L32").Value = 29030
$ws.Range("M32").Value = -2560.6924
$ws.Range("N32").Value = -29604
# Row 88
$ws.Range("H88").Value = 62502110
$ws.Range("I88").Value = 1633.8889
$ws.Range("J88").Value = 142859870
$ws.Range("K88").Value = 1633.8889
$ws.Range("L88").Value = 142859870
$ws.Range("M88").Value = -1227.8889
$ws.Range("N88").Value = -142860682
# Row 91
$ws.Range("H91").Value = 62502110
$ws.Range("I91").Value = 1633.8889
$ws.Range("J91").Value = 142859870
$ws.Range("K91").Value = 1633.8889
$ws.Range("L91").Value = 142859870
$ws.Range("M91").Value = -229.8888999999999
$ws.Range("N91").Value = -142862678
# Row 110
$ws.Range("H110").Value = 7636.263
$ws.Range("I110").Value = 9136.429
$ws.Range("K110").Value = 9136.429
$ws.Range("M110").Value = -7091.429

$ws = $wb.Worksheets.Item("BSM")
# Row 15
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
# Row 86
$ws.Range("H86").Value = 33335650
$ws.Range("I86").Value = 40001680
$ws.Range("J86").Value = 5500
$ws.Range("K86").Value = 40001680
$ws.Range("L86").Value = 5500
$ws.Range("M86").Value = -40000557
$ws.Range("N86").Value = -7746
# Row 89
$ws.Range("H89").Value = 33335650
$ws.Range("I89").Value = 40001680
$ws.Range("J89").Value = 5500
$ws.Range("K89").Value = 200008400
$ws.Range("L89").Value = 27500
$ws.Range("M89").Value = -200002784
$ws.Range("N89").Value = -38732

$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 90911860
$ws.Range("I62").Value = 2984.1667
$ws.Range("J62").Value = 200002500
$ws.Range("K62").Value = 2984.1667
$ws.Range("L62").Value = 200002500
$ws.Range("M62").Value = -2360.1667
$ws.Range("N62").Value = -200003748
# Row 65
$ws.Range("H65").Value = 90911860
$ws.Range("I65").Value = 2984.1667
$ws.Range("J65").Value = 200002500
$ws.Range("K65").Value = 14920.8335
$ws.Range("L65").Value = 1000012500
$ws.Range("M65").Value = -11800.8335
$ws.Range("N65").Value = -1000018740
# Row 107
$ws.Range("H107").Value = 29413426
$ws.Range("I107").Value = 38463020
$ws.Range("J107").Value = 2245
$ws.Range("K107").Value = 38463020
$ws.Range("L107").Value = 2245
$ws.Range("M107").Value = -38461100
$ws.Range("N107").Value = -6085

$ws = $wb.Worksheets.Item("CUL")
# Row 9
$ws.Range("H9").Value = 75574.94
$ws.Range("I9").Value = 223200.4
$ws.Range("J9").Value = 8472.454
$ws.Range("K9").Value = 669601.2
$ws.Range("L9").Value = 25417.362
$ws.Range("M9").Value = -669377.2
$ws.Range("N9").Value = -25865.362
# Row 26
$ws.Range("H26").Value = 250.55556
$ws.Range("I26").Value = 210.16667
$ws.Range("J26").Value = 331.33334
$ws.Range("K26").Value = 630.50001
$ws.Range("L26").Value = 994.0000200000001
$ws.Range("M26").Value = -342.50001
$ws.Range("N26").Value = -1570.00002
# Row 33
$ws.Range("H33").Value = 4534.391
$ws.Range("I33").Value = 269.18182
$ws.Range("J33").Value = 8444.167
$ws.Range("K33").Value = 1615.09092
$ws.Range("L33").Value = 50665.00199999999
$ws.Range("M33").Value = -1332.09092
$ws.Range("N33").Value = -51231.00199999999
# Row 92
$ws.Range("H92").Value = 569.3333
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 569.3333
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 1707.9999
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -4203.9999
# Row 131
$ws.Range("H131").Value = 9260213
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 9260213
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 27780639
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -27790719

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 26671136
$ws.Range("I70").Value = 40004276
$ws.Range("J70").Value = 4853.9
$ws.Range("K70").Value = 40004276
$ws.Range("L70").Value = 4853.9
$ws.Range("M70").Value = -40004006
$ws.Range("N70").Value = -5393.9
# Row 73
$ws.Range("H73").Value = 26671136
$ws.Range("I73").Value = 40004276
$ws.Range("J73").Value = 4853.9
$ws.Range("K73").Value = 40004276
$ws.Range("L73").Value = 4853.9
$ws.Range("M73").Value = -40003340
$ws.Range("N73").Value = -6725.9
# Row 113
$ws.Range("H113").Value = 1668.5883
$ws.Range("I113").Value = 1207.8889
$ws.Range("J113").Value = 2186.875
$ws.Range("K113").Value = 1207.8889
$ws.Range("L113").Value = 2186.875
$ws.Range("M113").Value = 962.1111000000001
$ws.Range("N113").Value = -6526.875

$ws = $wb.Worksheets.Item("LTW")
# Row 136
$ws.Range("H136").Value = 2663380.5
$ws.Range("I136").Value = 5214402
$ws.Range("J136").Value = 1445
$ws.Range("K136").Value = 15643206
$ws.Range("L136").Value = 4335
$ws.Range("N136").Value = -9435

